$wb = $excel.ActiveWorkbook

# --- AddNewUser sheet changes (selection + D2 value) ---
$wsAddNewUser = $wb.Worksheets.Item("AddNewUser")
$wsAddNewUser.Range("D2").Value = 1289

# --- Insert three new worksheets after "AddNewUser" ---
$wsCase = $wb.Worksheets.Add($null, $wsAddNewUser)
$wsCase.Name = "AddNewCase"
$wsCase.Range("A1").Value = "AssignedTo"
$wsCase.Range("B1").Value = "CaseName"
$wsCase.Range("C1").Value = "Description"

$wsTemplate = $wb.Worksheets.Add($null, $wsCase)
$wsTemplate.Name = "AddNewH1Template"
$wsTemplate.Range("A1").Value = "TemplateName"
$wsTemplate.Range("B1").Value = "Description"

$wsPatrol = $wb.Worksheets.Add($null, $wsTemplate)
$wsPatrol.Name = "AddH1PatrolUnit"
$wsPatrol.Range("A1").Value = "PatrolUnitName"
$wsPatrol.Range("B1").Value = "Description"

# --- Fill in the data rows (second pass, matching original authoring order) ---
$wsCase.Range("A2").Value = 4777
$wsCase.Range("B2").Value = "testcase108"
$wsCase.Range("C2").Value = "testcase108"
$wsCase.Columns.Item(1).ColumnWidth = 10.3
$wsCase.Columns.Item(2).ColumnWidth = 10.5
$wsCase.Columns.Item(3).ColumnWidth = 10.5

$wsTemplate.Range("A2").Value = "test017"
$wsTemplate.Range("B2").Value = "test017"

$wsPatrol.Range("A2").Value = "test018"
$wsPatrol.Range("B2").Value = "test018"

# --- Set selections on each sheet (activating as needed) ---
[void]$wsAddNewUser.Activate()
$wsAddNewUser.Range("E8").Select() | Out-Null

[void]$wsCase.Activate()
$wsCase.Range("C8").Select() | Out-Null

[void]$wsTemplate.Activate()
$wsTemplate.Range("B2").Select() | Out-Null

[void]$wsPatrol.Activate()
$wsPatrol.Range("C8").Select() | Out-Null

# --- WebURL sheet keeps its own selection, just loses tab-selected status ---
$wsWebURL = $wb.Worksheets.Item("WebURL")
[void]$wsWebURL.Activate()
$wsWebURL.Range("D9").Select() | Out-Null

# --- Finally, PermissionUser becomes the active/selected tab ---
$wsPermissionUser = $wb.Worksheets.Item("PermissionUser")
[void]$wsPermissionUser.Activate()
